$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 863, shifting the existing rows 863:905 down by one.
$ws.Rows.Item(863).Insert()

# Populate the newly inserted row 863 with the new data point.
# Use a leading apostrophe so the date-like string is kept as literal text
# (matching the source file, where column A is plain text, not a real date),
# then strip the resulting quote-prefix/number formatting so the cell stays
# styled like its unstyled neighbours.
$ws.Cells.Item(863, 1).Value = "'2026/02/24"
$ws.Cells.Item(863, 1).ClearFormats()
$ws.Cells.Item(863, 2).Value = "火"
$ws.Cells.Item(863, 3).Value = 12
$ws.Cells.Item(863, 4).Value = 201
